$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Email value (B4): replace the old test email with the new one.
# Clear its existing (non-visual) format first so the new hyperlink format
# isn't polluted by the old "applyFill" flag.
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "name.test12@gmail.com"

# The Password cell (B6) used to share the same (invisible) fill-applying
# style as B4; reset it back to the plain "Normal" style now that B4 has
# moved on to its own (hyperlink) style.
$ws.Range("B6").Style = "Normal"

# Typing an e-mail address into a cell makes Excel auto-create a mailto:
# hyperlink (new Hyperlink cell style/font + a sheet-level hyperlink entry).
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:name.test12@gmail.com")

# The active selection ends up on B4 after the edit.
$ws.Range("B4").Select()
